$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell in the workbook that shows the old status gets the new one)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate "Latest Target File" (F) / "Latest Handback File" (G) columns
#    for the handed-back rows, and stamp "Latest Handback DateTime" (H).
# ---------------------------------------------------------------------------

$mdFileName = "597a00e2-acda-459d-b80a-0c78127569f3.md"
$zhXlfFileName = "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.zh-cn.xlf"
$deXlfFileName = "597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/31a0c34aeee7c85a42f71fd494127562a9aa5204/e2e/597a00e2-acda-459d-b80a-0c78127569f3.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b4951ad2042542b9db9179ea219fa130609cc4e3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca8425e1e9c32ce3e5bea1ab186b39126c3afdb1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/597a00e2-acda-459d-b80a-0c78127569f3.c0e5cbbd5c3502ce52a22c5403ff31242e293cbc.de-de.xlf"

# -- zh-cn sheet --------------------------------------------------------
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl, "", "", $mdFileName)
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", $zhXlfFileName)
$wsZh.Range("G2").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl, "", "", $mdFileName)
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", $zhXlfFileName)
$wsZh.Range("G3").Style = "HyperLink"

$wsZh.Range("H2").Value = "2016-03-20 18:50:21"
$wsZh.Range("H3").Value = "2016-03-20 18:50:21"

# -- de-de sheet --------------------------------------------------------
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl, "", "", $mdFileName)
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", $deXlfFileName)
$wsDe.Range("G2").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl, "", "", $mdFileName)
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", $deXlfFileName)
$wsDe.Range("G3").Style = "HyperLink"

$wsDe.Range("H2").Value = "2016-03-20 18:50:27"
$wsDe.Range("H3").Value = "2016-03-20 18:50:27"
